$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 264
$ws.Range("F4").Value = 863
$ws.Range("F5").Value = 243
$ws.Range("F7").Value = 638
$ws.Range("F8").Value = 231
$ws.Range("F10").Value = 368
$ws.Range("F11").Value = 169
$ws.Range("F12").Value = 741
$ws.Range("F13").Value = 100
$ws.Range("F14").Value = 1873
$ws.Range("F15").Value = 392
$ws.Range("F16").Value = 4641
$ws.Range("F17").Value = 396
$ws.Range("F18").Value = 501
$ws.Range("F19").Value = 23
$ws.Range("F21").Value = 156

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 104
$ws.Range("F14").Value = 44
$ws.Range("F17").Value = 32

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5399
$ws.Range("F3").Value = 345
$ws.Range("F4").Value = 316

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5399
$ws.Range("F4").Value = 345
$ws.Range("F6").Value = 316
$ws.Range("F7").Value = 264
$ws.Range("F13").Value = 863
$ws.Range("F16").Value = 243
$ws.Range("F18").Value = 638
$ws.Range("F19").Value = 231
$ws.Range("F22").Value = 368
$ws.Range("F23").Value = 169
$ws.Range("F26").Value = 741
$ws.Range("F27").Value = 100
$ws.Range("F28").Value = 104
$ws.Range("F29").Value = 1873
$ws.Range("F30").Value = 392
$ws.Range("F31").Value = 4641
$ws.Range("F32").Value = 44
$ws.Range("F33").Value = 396
$ws.Range("F34").Value = 501
$ws.Range("F35").Value = 23
$ws.Range("F38").Value = 156
$ws.Range("F41").Value = 32
